# Applies the "spell/grammar check pass" diff to progress.docx:
# text content is unchanged, but several runs get split and
# <w:proofErr/> spellStart/spellEnd/gramStart/gramEnd markers are
# inserted around individual words, exactly as Word's background
# proofer does when a document is opened/edited and re-saved.

$d = $word.ActiveDocument

# ---------------------------------------------------------------
# Helper: replace the *content* (runs only - not the enclosing
# <w:p>, so paragraph properties such as numPr/pStyle survive) of
# the range [start,end) with an explicit run of WordprocessingML.
# $innerXml is everything that belongs *inside* <w:p> ... </w:p>.
# ---------------------------------------------------------------
function Set-RunXml([int]$start, [int]$end, [string]$innerXml) {
    $rng = $d.Range($start, $end)
    $frag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p>' + $innerXml + '</w:p></w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($frag)
}

function NeedsPreserve([string]$text) {
    if ($text.Length -eq 0) { return $true }
    if ($text.Substring(0,1) -eq " " -or $text.Substring($text.Length-1,1) -eq " ") { return $true }
    return $false
}

function TEl([string]$text) {
    if (NeedsPreserve $text) {
        return '<w:t xml:space="preserve">' + $text + '</w:t>'
    }
    return '<w:t>' + $text + '</w:t>'
}

function R([string]$text) {
    return '<w:r>' + $(TEl $text) + '</w:r>'
}

function RLang([string]$text) {
    return '<w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr>' + $(TEl $text) + '</w:r>'
}

$SpellStart = '<w:proofErr w:type="spellStart"/>'
$SpellEnd   = '<w:proofErr w:type="spellEnd"/>'
$GramStart  = '<w:proofErr w:type="gramStart"/>'
$GramEnd    = '<w:proofErr w:type="gramEnd"/>'
$Dash = [string][char]0x2013

# 1. "Findings:" -> "Findings" (spellStart/spellEnd) + ":"
$r1a = RLang "Findings"
$r1b = RLang ":"
$inner1 = $SpellStart + $r1a + $SpellEnd + $r1b
Set-RunXml 0 9 $inner1

# 3. "...128x128. So I stick..." -> split out "So" as a grammar error
$r3a = R "64x64 does not result in any noticeable lower accuracy over 128x128. "
$r3b = R "So"
$r3c = R " I stick with 64x64 for a fair comparison with DCGAN. "
$inner3 = $r3a + $GramStart + $r3b + $GramEnd + $r3c
Set-RunXml 107 232 $inner3

# 5. "ReLU " / "improve true positive rate over LeakyReLU"
$r5a = R "ReLU"
$r5b = R " "
$r5c = R "improve true positive rate over "
$r5d = R "LeakyReLU"
$inner5 = $SpellStart + $r5a + $SpellEnd + $r5b + $r5c + $SpellStart + $r5d + $SpellEnd
Set-RunXml 350 396 $inner5

# 9. "001: uses rgb and a large batch size of 200, activation = relu."
$r9a = R "001: uses "
$r9b = R "rgb"
$r9c = R " and a large b"
$r9d = R "atch size of 200"
$r9e = R ", activation = "
$r9f = R "relu"
$r9g = R "."
$inner9 = $r9a + $SpellStart + $r9b + $SpellEnd + $r9c + $r9d + $r9e + $SpellStart + $r9f + $SpellEnd + $r9g
Set-RunXml 558 621 $inner9

# 10. "002: uses grayscale and a batch size of 25, activation = relu. "
$r10a = R "002: uses grayscale and a batch size of 25"
$r10b = R ", activation = "
$r10c = R "relu"
$r10d = R ". "
$inner10 = $r10a + $r10b + $SpellStart + $r10c + $SpellEnd + $r10d
Set-RunXml 622 685 $inner10

# 11. "003: uses grayscale and a batch size of 25, activation = relu.   "
$r11a = R "003: uses grayscale and a batch size of 25"
$r11b = R ", activation = "
$r11c = R "relu"
$r11d = R ". "
$r11e = R " "
$inner11 = $r11a + $r11b + $SpellStart + $r11c + $SpellEnd + $r11d + $r11e
Set-RunXml 686 750 $inner11

# 13. "Grayscale, batch size=25, activation=relu - v002,003,004: "
$r13a = R "Grayscale, batch size=25, activation="
$r13b = R "relu"
$r13tail = " " + $Dash + " v002,003,004: "
$r13c = R $r13tail
$inner13 = $r13a + $SpellStart + $r13b + $SpellEnd + $r13c
Set-RunXml 762 820 $inner13

# 14. "Multi-Acc: 9525-9530"
$r14a = R "Multi-"
$r14b = R "Acc"
$r14c = R ": 9525-9530"
$inner14 = $r14a + $SpellStart + $r14b + $SpellEnd + $r14c
Set-RunXml 821 841 $inner14

# 16. "Acc: 9812-9838"
$r16a = R "Acc"
$r16b = R ": 9812-9838"
$inner16 = $SpellStart + $r16a + $SpellEnd + $r16b
Set-RunXml 857 871 $inner16

# 19. "Grayscale, batch size=25, activation=leakyrelu - v005,006,007: "
$r19a = R "Grayscale, batch size=25, activation="
$r19b = R "leakyrelu"
$r19tail = " " + $Dash + " v005,006,007: "
$r19c = R $r19tail
$inner19 = $r19a + $SpellStart + $r19b + $SpellEnd + $r19c
Set-RunXml 901 964 $inner19

# 20. "Multi-Acc: 9493-9503"
$r20a = R "Multi-"
$r20b = R "Acc"
$r20c = R ":"
$r20d = R " 9493-9503"
$inner20 = $r20a + $SpellStart + $r20b + $SpellEnd + $r20c + $r20d
Set-RunXml 965 985 $inner20

# 22. "Acc: 9809-9838"
$r22a = R "Acc"
$r22b = R ": "
$r22c = R "9809-9838"
$inner22 = $SpellStart + $r22a + $SpellEnd + $r22b + $r22c
Set-RunXml 1001 1015 $inner22

Write-Output "done"
